$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.7321483333333333
$ws.Range("H2").Value = 2.196445
$ws.Range("I2").Value = 0.05113520435363902
$ws.Range("J2").Value = 0.05113520435363902
$ws.Range("M2").Value = 28.31444233333334
$ws.Range("N2").Value = 84.94332700000001
$ws.Range("O2").Value = 0.2747173016130739
$ws.Range("P2").Value = 0.2747173016130739
$ws.Range("Q2").Value = 20.73037176361278
$ws.Range("R2").Value = 186.573345872515
$ws.Range("S2").Value = 0.01404772535746482
$ws.Range("T2").Value = 0.01404772535746482
$ws.Range("G3").Value = 0.7321483333333333
$ws.Range("H3").Value = 2.196445
$ws.Range("I3").Value = 0.05113520435363902
$ws.Range("J3").Value = 0.05113520435363902
$ws.Range("O3").Value = 0.2090339131726295
$ws.Range("P3").Value = 0.2090339131726295
$ws.Range("Q3").Value = 15.77385445265722
$ws.Range("R3").Value = 141.964690073915
$ws.Range("S3").Value = 0.01068899186692325
$ws.Range("T3").Value = 0.01068899186692325
$ws.Range("G4").Value = 0.7321483333333333
$ws.Range("H4").Value = 2.196445
$ws.Range("I4").Value = 0.05113520435363902
$ws.Range("J4").Value = 0.05113520435363902
$ws.Range("M4").Value = 5.413469333333334
$ws.Range("N4").Value = 16.240408
$ws.Range("O4").Value = 0.0525235026743817
$ws.Range("P4").Value = 0.0525235026743817
$ws.Range("Q4").Value = 3.963462549951112
$ws.Range("R4").Value = 35.67116294956
$ws.Range("S4").Value = 0.002685800042623414
$ws.Range("T4").Value = 0.002685800042623414
$ws.Range("G5").Value = 0.7321483333333333
$ws.Range("H5").Value = 2.196445
$ws.Range("I5").Value = 0.05113520435363902
$ws.Range("J5").Value = 0.05113520435363902
$ws.Range("M5").Value = 47.79503400000001
$ws.Range("N5").Value = 143.385102
$ws.Range("O5").Value = 0.4637252825399149
$ws.Range("P5").Value = 0.4637252825399149
$ws.Range("Q5").Value = 34.99305448471
$ws.Range("R5").Value = 314.93749036239
$ws.Range("S5").Value = 0.02371268708662754
$ws.Range("T5").Value = 0.02371268708662754
$ws.Range("I6").Value = 0.7165747117895102
$ws.Range("J6").Value = 0.7165747117895102
$ws.Range("M6").Value = 28.31444233333334
$ws.Range("N6").Value = 84.94332700000001
$ws.Range("O6").Value = 0.2747173016130739
$ws.Range("P6").Value = 0.2747173016130739
$ws.Range("Q6").Value = 290.5016291529318
$ws.Range("R6").Value = 2614.514662376386
$ws.Range("S6").Value = 0.1968554712269804
$ws.Range("T6").Value = 0.1968554712269804
$ws.Range("I7").Value = 0.7165747117895102
$ws.Range("J7").Value = 0.7165747117895102
$ws.Range("O7").Value = 0.2090339131726295
$ws.Range("P7").Value = 0.2090339131726295
$ws.Range("S7").Value = 0.1497884160859105
$ws.Range("T7").Value = 0.1497884160859105
$ws.Range("I8").Value = 0.7165747117895102
$ws.Range("J8").Value = 0.7165747117895102
$ws.Range("M8").Value = 5.413469333333334
$ws.Range("N8").Value = 16.240408
$ws.Range("O8").Value = 0.0525235026743817
$ws.Range("P8").Value = 0.0525235026743817
$ws.Range("Q8").Value = 55.54132559592713
$ws.Range("R8").Value = 499.8719303633441
$ws.Range("S8").Value = 0.03763701379107064
$ws.Range("T8").Value = 0.03763701379107064
$ws.Range("I9").Value = 0.7165747117895102
$ws.Range("J9").Value = 0.7165747117895102
$ws.Range("M9").Value = 47.79503400000001
$ws.Range("N9").Value = 143.385102
$ws.Range("O9").Value = 0.4637252825399149
$ws.Range("P9").Value = 0.4637252825399149
$ws.Range("Q9").Value = 490.3693697712041
$ws.Range("R9").Value = 4413.324327940836
$ws.Range("S9").Value = 0.3322938106855486
$ws.Range("T9").Value = 0.3322938106855486
$ws.Range("G10").Value = 2.568000333333333
$ws.Range("H10").Value = 7.704001
$ws.Range("I10").Value = 0.1793560346266988
$ws.Range("J10").Value = 0.1793560346266988
$ws.Range("M10").Value = 28.31444233333334
$ws.Range("N10").Value = 84.94332700000001
$ws.Range("O10").Value = 0.2747173016130739
$ws.Range("P10").Value = 0.2747173016130739
$ws.Range("Q10").Value = 72.71149735014745
$ws.Range("R10").Value = 654.403476151327
$ws.Range("S10").Value = 0.04927220586066774
$ws.Range("T10").Value = 0.04927220586066774
$ws.Range("G11").Value = 2.568000333333333
$ws.Range("H11").Value = 7.704001
$ws.Range("I11").Value = 0.1793560346266988
$ws.Range("J11").Value = 0.1793560346266988
$ws.Range("O11").Value = 0.2090339131726295
$ws.Range("P11").Value = 0.2090339131726295
$ws.Range("Q11").Value = 55.32658021353856
$ws.Range("R11").Value = 497.939221921847
$ws.Range("S11").Value = 0.03749149376914448
$ws.Range("T11").Value = 0.03749149376914448
$ws.Range("G12").Value = 2.568000333333333
$ws.Range("H12").Value = 7.704001
$ws.Range("I12").Value = 0.1793560346266988
$ws.Range("J12").Value = 0.1793560346266988
$ws.Range("M12").Value = 5.413469333333334
$ws.Range("N12").Value = 16.240408
$ws.Range("O12").Value = 0.0525235026743817
$ws.Range("P12").Value = 0.0525235026743817
$ws.Range("Q12").Value = 13.90179105248978
$ws.Range("R12").Value = 125.116119472408
$ws.Range("S12").Value = 0.00942040716438191
$ws.Range("T12").Value = 0.00942040716438191
$ws.Range("G13").Value = 2.568000333333333
$ws.Range("H13").Value = 7.704001
$ws.Range("I13").Value = 0.1793560346266988
$ws.Range("J13").Value = 0.1793560346266988
$ws.Range("M13").Value = 47.79503400000001
$ws.Range("N13").Value = 143.385102
$ws.Range("O13").Value = 0.4637252825399149
$ws.Range("P13").Value = 0.4637252825399149
$ws.Range("Q13").Value = 122.737663243678
$ws.Range("R13").Value = 1104.638969193102
$ws.Range("S13").Value = 0.08317192783250464
$ws.Range("T13").Value = 0.08317192783250464
$ws.Range("G14").Value = 0.7579039999999999
$ws.Range("H14").Value = 2.273712
$ws.Range("I14").Value = 0.05293404923015203
$ws.Range("J14").Value = 0.05293404923015203
$ws.Range("M14").Value = 28.31444233333334
$ws.Range("N14").Value = 84.94332700000001
$ws.Range("O14").Value = 0.2747173016130739
$ws.Range("P14").Value = 0.2747173016130739
$ws.Range("Q14").Value = 21.45962910220267
$ws.Range("R14").Value = 193.136661919824
$ws.Range("S14").Value = 0.01454189916796098
$ws.Range("T14").Value = 0.01454189916796098
$ws.Range("G15").Value = 0.7579039999999999
$ws.Range("H15").Value = 2.273712
$ws.Range("I15").Value = 0.05293404923015203
$ws.Range("J15").Value = 0.05293404923015203
$ws.Range("O15").Value = 0.2090339131726295
$ws.Range("P15").Value = 0.2090339131726295
$ws.Range("Q15").Value = 16.32875039222933
$ws.Range("R15").Value = 146.958753530064
$ws.Range("S15").Value = 0.0110650114506513
$ws.Range("T15").Value = 0.0110650114506513
$ws.Range("G16").Value = 0.7579039999999999
$ws.Range("H16").Value = 2.273712
$ws.Range("I16").Value = 0.05293404923015203
$ws.Range("J16").Value = 0.05293404923015203
$ws.Range("M16").Value = 5.413469333333334
$ws.Range("N16").Value = 16.240408
$ws.Range("O16").Value = 0.0525235026743817
$ws.Range("P16").Value = 0.0525235026743817
$ws.Range("Q16").Value = 4.102890061610667
$ws.Range("R16").Value = 36.926010554496
$ws.Range("S16").Value = 0.002780281676305743
$ws.Range("T16").Value = 0.002780281676305743
$ws.Range("G17").Value = 0.7579039999999999
$ws.Range("H17").Value = 2.273712
$ws.Range("I17").Value = 0.05293404923015203
$ws.Range("J17").Value = 0.05293404923015203
$ws.Range("M17").Value = 47.79503400000001
$ws.Range("N17").Value = 143.385102
$ws.Range("O17").Value = 0.4637252825399149
$ws.Range("P17").Value = 0.4637252825399149
$ws.Range("Q17").Value = 122.737663243678
$ws.Range("R17").Value = 1104.638969193102
$ws.Range("S17").Value = 0.08317192783250464
$ws.Range("T17").Value = 0.08317192783250464
